$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The roster table (A2:C18) was refreshed with a new snapshot of player data.
# Column A (player name) keeps referencing the same shared-string slot order
# as before, but the underlying text for each row changes together with
# columns B (position) and C (team) below.

$ws.Range("A2").Value = "Jordan Poole"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Washington Wizards"

$ws.Range("A3").Value = "Kyrie Irving"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "Dallas Mavericks"

$ws.Range("A4").Value = "Malcolm Brogdon"
$ws.Range("B4").Value = "PG,SG"
$ws.Range("C4").Value = "Detroit Pistons"

$ws.Range("A5").Value = "Tobias Harris"
$ws.Range("B5").Value = "SF,PF"
$ws.Range("C5").Value = "Detroit Pistons"

$ws.Range("A6").Value = "John Collins"
$ws.Range("B6").Value = "PF,C"
$ws.Range("C6").Value = "Utah Jazz"

$ws.Range("A7").Value = "Lauri Markkanen"
$ws.Range("B7").Value = "SF,PF"
$ws.Range("C7").Value = "Utah Jazz"

$ws.Range("A8").Value = "Jonas Valanciunas"
$ws.Range("B8").Value = "C"
$ws.Range("C8").Value = "Washington Wizards"

$ws.Range("A9").Value = "Keyonte George"
$ws.Range("B9").Value = "PG,SG"
$ws.Range("C9").Value = "Utah Jazz"

$ws.Range("A10").Value = "Shai Gilgeous-Alexander"
$ws.Range("B10").Value = "PG"
$ws.Range("C10").Value = "Oklahoma City Thunder"

$ws.Range("A11").Value = "RJ Barrett"
$ws.Range("B11").Value = "SF,PF"
$ws.Range("C11").Value = "Toronto Raptors"

$ws.Range("A12").Value = "Jimmy Butler"
$ws.Range("B12").Value = "SF,PF"
$ws.Range("C12").Value = "Miami Heat"

$ws.Range("A13").Value = "Zach LaVine"
$ws.Range("B13").Value = "SG,SF"
$ws.Range("C13").Value = "Chicago Bulls"

$ws.Range("A14").Value = "CJ McCollum"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "New Orleans Pelicans"

$ws.Range("A15").Value = "Jalen Williams"
$ws.Range("B15").Value = "SG,SF,PF,C"
$ws.Range("C15").Value = "Oklahoma City Thunder"

$ws.Range("A16").Value = "Christian Braun"
$ws.Range("B16").Value = "SG,SF"
$ws.Range("C16").Value = "Denver Nuggets"

$ws.Range("A17").Value = "Joel Embiid"
$ws.Range("B17").Value = "C"
$ws.Range("C17").Value = "Philadelphia 76ers"

$ws.Range("A18").Value = "Dennis Schröder"
$ws.Range("B18").Value = "PG"
$ws.Range("C18").Value = "Brooklyn Nets"
